# Add a new "spare_parts" worksheet right after "tasks" and before
# "task_compatibility", populate it with the SET header + S1..S4 members,
# and make it the active/selected sheet (matching the authored diff).

$wb = $excel.ActiveWorkbook

$tasksSheet = $wb.Worksheets.Item("tasks")

# Insert the new sheet immediately after "tasks"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tasksSheet)
$newSheet.Name = "spare_parts"

$newSheet.Range("A1").Value = "SET"
$newSheet.Range("A2").Value = "S1"
$newSheet.Range("A3").Value = "S2"
$newSheet.Range("A4").Value = "S3"
$newSheet.Range("A5").Value = "S4"

# Make the new sheet the active tab / selected tab, with B1 selected,
# matching the target workbook state.
$newSheet.Activate()
$newSheet.Range("B1").Select()
